# Generate Report for Handoff
#
# The report rotates: the row that used to describe the "cf71a493" file
# (which had just been handed back) now drops to the bottom as a fresh
# "Ready for handoff" entry with new handoff timestamps, while the two
# rows that used to describe "ffff652c5f9a" / "fffffff5fc5ce1" shift up
# one slot (their data is unchanged, only which row shows which file).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "ffff652c5f9a-3d0e-4bdf-99c0-9656fdd39890.md"
$ov.Range("B2").Value = "Handed back: in sync with en-US"
$ov.Range("C2").Value = "Handed back: in sync with en-US"
$ov.Range("D2").Value = "2016-03-23 23:10:57"

$ov.Range("A3").Value = "fffffff5fc5ce1-27ae-424d-b80e-a1ad0cd9e8af.md"
$ov.Range("B3").Value = "Handed back: in sync with en-US"
$ov.Range("C3").Value = "Handed back: in sync with en-US"
$ov.Range("D3").Value = "2016-03-23 23:10:57"

$ov.Range("A4").Value = "cf71a493-e1f2-42f7-a347-cd876ddb551c.md"
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"
$ov.Range("D4").Value = "2016-03-23 23:14:21"

$ov.Hyperlinks.Item(1).TextToDisplay = "ffff652c5f9a-3d0e-4bdf-99c0-9656fdd39890.md"
$ov.Hyperlinks.Item(2).TextToDisplay = "fffffff5fc5ce1-27ae-424d-b80e-a1ad0cd9e8af.md"
$ov.Hyperlinks.Item(3).TextToDisplay = "cf71a493-e1f2-42f7-a347-cd876ddb551c.md"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "ffff652c5f9a-3d0e-4bdf-99c0-9656fdd39890.md"
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("D2").Value = "e498de5f-2453-4d31-81c0-8ce5f05474c0.88694f32d0f0406db45b4acd55aaf8271b2f35be.zh-cn.xlf"
$zh.Range("E2").Value = "2016-03-23 23:10:52"
$zh.Range("F2").Value = "e498de5f-2453-4d31-81c0-8ce5f05474c0.md"
$zh.Range("G2").Value = "e498de5f-2453-4d31-81c0-8ce5f05474c0.88694f32d0f0406db45b4acd55aaf8271b2f35be.zh-cn.xlf"
$zh.Range("H2").Value = "2016-03-23 23:11:24"
$zh.Range("J2").Value = "Include"

$zh.Range("A3").Value = "fffffff5fc5ce1-27ae-424d-b80e-a1ad0cd9e8af.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Handed back: in sync with en-US"
$zh.Range("D3").Value = "e498de5f-2453-4d31-81c0-8ce5f05474c0.88694f32d0f0406db45b4acd55aaf8271b2f35be.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-23 23:10:52"
$zh.Range("F3").Value = "e498de5f-2453-4d31-81c0-8ce5f05474c0.md"
$zh.Range("G3").Value = "e498de5f-2453-4d31-81c0-8ce5f05474c0.88694f32d0f0406db45b4acd55aaf8271b2f35be.zh-cn.xlf"
$zh.Range("H3").Value = "2016-03-23 23:11:24"
$zh.Range("J3").Value = "Include"

$zh.Range("A4").Value = "cf71a493-e1f2-42f7-a347-cd876ddb551c.md"
$zh.Range("B4").Value = ".md"
$zh.Range("C4").Value = "Ready for handoff"
$zh.Range("D4").Value = "cf71a493-e1f2-42f7-a347-cd876ddb551c.992a5591db3e54946ed59a1de6fe4cbca95382ae.zh-cn.xlf"
$zh.Range("E4").Value = "2016-03-23 23:14:12"
$zh.Range("F4").Value = "cf71a493-e1f2-42f7-a347-cd876ddb551c.md"
$zh.Range("G4").Value = "cf71a493-e1f2-42f7-a347-cd876ddb551c.992a5591db3e54946ed59a1de6fe4cbca95382ae.zh-cn.xlf"
$zh.Range("H4").Value = "2016-03-23 23:13:07"
$zh.Range("J4").Value = "Include"

$zh.Hyperlinks.Item(1).TextToDisplay = "ffff652c5f9a-3d0e-4bdf-99c0-9656fdd39890.md"
$zh.Hyperlinks.Item(2).TextToDisplay = "e498de5f-2453-4d31-81c0-8ce5f05474c0.88694f32d0f0406db45b4acd55aaf8271b2f35be.zh-cn.xlf"
$zh.Hyperlinks.Item(3).TextToDisplay = "e498de5f-2453-4d31-81c0-8ce5f05474c0.md"
$zh.Hyperlinks.Item(4).TextToDisplay = "e498de5f-2453-4d31-81c0-8ce5f05474c0.88694f32d0f0406db45b4acd55aaf8271b2f35be.zh-cn.xlf"
$zh.Hyperlinks.Item(5).TextToDisplay = "fffffff5fc5ce1-27ae-424d-b80e-a1ad0cd9e8af.md"
$zh.Hyperlinks.Item(6).TextToDisplay = "e498de5f-2453-4d31-81c0-8ce5f05474c0.88694f32d0f0406db45b4acd55aaf8271b2f35be.zh-cn.xlf"
$zh.Hyperlinks.Item(7).TextToDisplay = "e498de5f-2453-4d31-81c0-8ce5f05474c0.md"
$zh.Hyperlinks.Item(8).TextToDisplay = "e498de5f-2453-4d31-81c0-8ce5f05474c0.88694f32d0f0406db45b4acd55aaf8271b2f35be.zh-cn.xlf"
$zh.Hyperlinks.Item(9).TextToDisplay = "cf71a493-e1f2-42f7-a347-cd876ddb551c.md"
$zh.Hyperlinks.Item(10).TextToDisplay = "cf71a493-e1f2-42f7-a347-cd876ddb551c.992a5591db3e54946ed59a1de6fe4cbca95382ae.zh-cn.xlf"
$zh.Hyperlinks.Item(11).TextToDisplay = "cf71a493-e1f2-42f7-a347-cd876ddb551c.md"
$zh.Hyperlinks.Item(12).TextToDisplay = "cf71a493-e1f2-42f7-a347-cd876ddb551c.992a5591db3e54946ed59a1de6fe4cbca95382ae.zh-cn.xlf"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "ffff652c5f9a-3d0e-4bdf-99c0-9656fdd39890.md"
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("D2").Value = "e498de5f-2453-4d31-81c0-8ce5f05474c0.88694f32d0f0406db45b4acd55aaf8271b2f35be.de-de.xlf"
$de.Range("E2").Value = "2016-03-23 23:10:57"
$de.Range("F2").Value = "e498de5f-2453-4d31-81c0-8ce5f05474c0.md"
$de.Range("G2").Value = "e498de5f-2453-4d31-81c0-8ce5f05474c0.88694f32d0f0406db45b4acd55aaf8271b2f35be.de-de.xlf"
$de.Range("H2").Value = "2016-03-23 23:11:31"
$de.Range("J2").Value = "Include"

$de.Range("A3").Value = "fffffff5fc5ce1-27ae-424d-b80e-a1ad0cd9e8af.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Handed back: in sync with en-US"
$de.Range("D3").Value = "e498de5f-2453-4d31-81c0-8ce5f05474c0.88694f32d0f0406db45b4acd55aaf8271b2f35be.de-de.xlf"
$de.Range("E3").Value = "2016-03-23 23:10:57"
$de.Range("F3").Value = "e498de5f-2453-4d31-81c0-8ce5f05474c0.md"
$de.Range("G3").Value = "e498de5f-2453-4d31-81c0-8ce5f05474c0.88694f32d0f0406db45b4acd55aaf8271b2f35be.de-de.xlf"
$de.Range("H3").Value = "2016-03-23 23:11:31"
$de.Range("J3").Value = "Include"

$de.Range("A4").Value = "cf71a493-e1f2-42f7-a347-cd876ddb551c.md"
$de.Range("B4").Value = ".md"
$de.Range("C4").Value = "Ready for handoff"
$de.Range("D4").Value = "cf71a493-e1f2-42f7-a347-cd876ddb551c.992a5591db3e54946ed59a1de6fe4cbca95382ae.de-de.xlf"
$de.Range("E4").Value = "2016-03-23 23:14:21"
$de.Range("F4").Value = "cf71a493-e1f2-42f7-a347-cd876ddb551c.md"
$de.Range("G4").Value = "cf71a493-e1f2-42f7-a347-cd876ddb551c.992a5591db3e54946ed59a1de6fe4cbca95382ae.de-de.xlf"
$de.Range("H4").Value = "2016-03-23 23:13:13"
$de.Range("J4").Value = "Include"

$de.Hyperlinks.Item(1).TextToDisplay = "ffff652c5f9a-3d0e-4bdf-99c0-9656fdd39890.md"
$de.Hyperlinks.Item(2).TextToDisplay = "e498de5f-2453-4d31-81c0-8ce5f05474c0.88694f32d0f0406db45b4acd55aaf8271b2f35be.de-de.xlf"
$de.Hyperlinks.Item(3).TextToDisplay = "e498de5f-2453-4d31-81c0-8ce5f05474c0.md"
$de.Hyperlinks.Item(4).TextToDisplay = "e498de5f-2453-4d31-81c0-8ce5f05474c0.88694f32d0f0406db45b4acd55aaf8271b2f35be.de-de.xlf"
$de.Hyperlinks.Item(5).TextToDisplay = "fffffff5fc5ce1-27ae-424d-b80e-a1ad0cd9e8af.md"
$de.Hyperlinks.Item(6).TextToDisplay = "e498de5f-2453-4d31-81c0-8ce5f05474c0.88694f32d0f0406db45b4acd55aaf8271b2f35be.de-de.xlf"
$de.Hyperlinks.Item(7).TextToDisplay = "e498de5f-2453-4d31-81c0-8ce5f05474c0.md"
$de.Hyperlinks.Item(8).TextToDisplay = "e498de5f-2453-4d31-81c0-8ce5f05474c0.88694f32d0f0406db45b4acd55aaf8271b2f35be.de-de.xlf"
$de.Hyperlinks.Item(9).TextToDisplay = "cf71a493-e1f2-42f7-a347-cd876ddb551c.md"
$de.Hyperlinks.Item(10).TextToDisplay = "cf71a493-e1f2-42f7-a347-cd876ddb551c.992a5591db3e54946ed59a1de6fe4cbca95382ae.de-de.xlf"
$de.Hyperlinks.Item(11).TextToDisplay = "cf71a493-e1f2-42f7-a347-cd876ddb551c.md"
$de.Hyperlinks.Item(12).TextToDisplay = "cf71a493-e1f2-42f7-a347-cd876ddb551c.992a5591db3e54946ed59a1de6fe4cbca95382ae.de-de.xlf"
